$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial of 45875 (2025-08-06) for rows
# 2 through 43; the update bumps it to 45877 (2025-08-08) for every row.
for ($r = 2; $r -le 43; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45875) {
        $cell.Value2 = 45877
    }
}
